$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 26).Value = "2025-10-17T07:09:39.947836"
$ws.Cells.Item(3, 26).Value = "2025-10-17T07:09:39.948840"
$ws.Cells.Item(4, 26).Value = "2025-10-17T07:09:39.948840"
$ws.Cells.Item(5, 26).Value = "2025-10-17T07:09:39.948840"
$ws.Cells.Item(6, 26).Value = "2025-10-17T07:09:39.948840"
$ws.Cells.Item(7, 26).Value = "2025-10-17T07:09:39.948840"
$ws.Cells.Item(8, 26).Value = "2025-10-17T07:09:39.949839"
$ws.Cells.Item(9, 26).Value = "2025-10-17T07:09:39.949839"
$ws.Cells.Item(10, 26).Value = "2025-10-17T07:09:39.949839"
$ws.Cells.Item(11, 26).Value = "2025-10-17T07:09:39.949839"
$ws.Cells.Item(12, 26).Value = "2025-10-17T07:09:39.950838"
$ws.Cells.Item(13, 26).Value = "2025-10-17T07:09:39.950838"
$ws.Cells.Item(14, 26).Value = "2025-10-17T07:09:39.950838"
$ws.Cells.Item(15, 26).Value = "2025-10-17T07:09:39.950838"
$ws.Cells.Item(16, 26).Value = "2025-10-17T07:09:39.950838"
$ws.Cells.Item(17, 26).Value = "2025-10-17T07:09:39.951841"
$ws.Cells.Item(18, 26).Value = "2025-10-17T07:09:39.951841"
$ws.Cells.Item(19, 26).Value = "2025-10-17T07:09:39.951841"
$ws.Cells.Item(20, 26).Value = "2025-10-17T07:09:39.952532"
$ws.Cells.Item(21, 26).Value = "2025-10-17T07:09:39.952532"
$ws.Cells.Item(22, 26).Value = "2025-10-17T07:09:39.953047"
$ws.Cells.Item(23, 26).Value = "2025-10-17T07:09:39.953047"
$ws.Cells.Item(24, 26).Value = "2025-10-17T07:09:39.953047"
$ws.Cells.Item(25, 26).Value = "2025-10-17T07:09:39.953047"
$ws.Cells.Item(26, 26).Value = "2025-10-17T07:09:39.953047"
$ws.Cells.Item(27, 26).Value = "2025-10-17T07:09:39.954066"
$ws.Cells.Item(28, 26).Value = "2025-10-17T07:09:39.954066"
$ws.Cells.Item(29, 26).Value = "2025-10-17T07:09:39.954066"
$ws.Cells.Item(30, 26).Value = "2025-10-17T07:09:39.954066"
$ws.Cells.Item(31, 26).Value = "2025-10-17T07:09:39.954066"
$ws.Cells.Item(32, 26).Value = "2025-10-17T07:09:39.955063"
$ws.Cells.Item(33, 26).Value = "2025-10-17T07:09:39.955063"
$ws.Cells.Item(34, 26).Value = "2025-10-17T07:09:39.955063"
$ws.Cells.Item(35, 26).Value = "2025-10-17T07:09:39.955063"
$ws.Cells.Item(36, 26).Value = "2025-10-17T07:09:39.955063"
$ws.Cells.Item(37, 26).Value = "2025-10-17T07:09:39.956062"
$ws.Cells.Item(38, 26).Value = "2025-10-17T07:09:39.956062"
$ws.Cells.Item(39, 26).Value = "2025-10-17T07:09:39.956062"
$ws.Cells.Item(40, 26).Value = "2025-10-17T07:09:39.956062"
$ws.Cells.Item(41, 26).Value = "2025-10-17T07:09:39.956812"
$ws.Cells.Item(42, 26).Value = "2025-10-17T07:09:39.956812"
$ws.Cells.Item(43, 26).Value = "2025-10-17T07:09:39.957486"
$ws.Cells.Item(44, 26).Value = "2025-10-17T07:09:39.957486"
$ws.Cells.Item(45, 26).Value = "2025-10-17T07:09:39.957486"
$ws.Cells.Item(46, 26).Value = "2025-10-17T07:09:40.036885"
$ws.Cells.Item(47, 26).Value = "2025-10-17T07:09:40.036885"
$ws.Cells.Item(48, 26).Value = "2025-10-17T07:09:40.036885"
$ws.Cells.Item(49, 26).Value = "2025-10-17T07:09:40.036885"
$ws.Cells.Item(50, 26).Value = "2025-10-17T07:09:40.037434"
$ws.Cells.Item(51, 26).Value = "2025-10-17T07:09:40.037434"
$ws.Cells.Item(52, 26).Value = "2025-10-17T07:09:40.037964"
$ws.Cells.Item(53, 26).Value = "2025-10-17T07:09:40.038105"
$ws.Cells.Item(54, 26).Value = "2025-10-17T07:09:40.038105"
$ws.Cells.Item(55, 26).Value = "2025-10-17T07:09:40.038105"
$ws.Cells.Item(56, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(57, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(58, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(59, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(60, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(61, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(62, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(63, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(64, 26).Value = "2025-10-17T07:09:40.038632"
$ws.Cells.Item(65, 26).Value = "2025-10-17T07:09:40.040618"
$ws.Cells.Item(66, 26).Value = "2025-10-17T07:09:40.040618"
$ws.Cells.Item(67, 26).Value = "2025-10-17T07:09:40.040618"
$ws.Cells.Item(68, 26).Value = "2025-10-17T07:09:40.040618"
$ws.Cells.Item(69, 26).Value = "2025-10-17T07:09:40.040618"
$ws.Cells.Item(70, 26).Value = "2025-10-17T07:09:40.041621"
$ws.Cells.Item(71, 26).Value = "2025-10-17T07:09:40.041621"
$ws.Cells.Item(72, 26).Value = "2025-10-17T07:09:40.041621"
$ws.Cells.Item(73, 26).Value = "2025-10-17T07:09:40.041621"
$ws.Cells.Item(74, 26).Value = "2025-10-17T07:09:40.041621"
$ws.Cells.Item(75, 26).Value = "2025-10-17T07:09:40.113817"
$ws.Cells.Item(76, 26).Value = "2025-10-17T07:09:40.116477"
$ws.Cells.Item(77, 26).Value = "2025-10-17T07:09:40.116477"
$ws.Cells.Item(78, 26).Value = "2025-10-17T07:09:40.117012"
$ws.Cells.Item(79, 26).Value = "2025-10-17T07:09:40.117012"
$ws.Cells.Item(80, 26).Value = "2025-10-17T07:09:40.117012"
$ws.Cells.Item(81, 26).Value = "2025-10-17T07:09:40.117012"
$ws.Cells.Item(82, 26).Value = "2025-10-17T07:09:40.117956"
$ws.Cells.Item(83, 26).Value = "2025-10-17T07:09:40.117956"
$ws.Cells.Item(84, 26).Value = "2025-10-17T07:09:40.117956"
$ws.Cells.Item(85, 26).Value = "2025-10-17T07:09:40.118502"
$ws.Cells.Item(86, 26).Value = "2025-10-17T07:09:40.118502"
$ws.Cells.Item(87, 26).Value = "2025-10-17T07:09:40.119013"
$ws.Cells.Item(88, 26).Value = "2025-10-17T07:09:40.119046"
$ws.Cells.Item(89, 26).Value = "2025-10-17T07:09:40.119046"
$ws.Cells.Item(90, 26).Value = "2025-10-17T07:09:40.119587"
$ws.Cells.Item(91, 26).Value = "2025-10-17T07:09:40.119587"
$ws.Cells.Item(92, 26).Value = "2025-10-17T07:09:40.119587"
$ws.Cells.Item(93, 26).Value = "2025-10-17T07:09:40.120129"
$ws.Cells.Item(94, 26).Value = "2025-10-17T07:09:40.120129"
$ws.Cells.Item(95, 26).Value = "2025-10-17T07:09:40.120129"
$ws.Cells.Item(96, 26).Value = "2025-10-17T07:09:40.120670"
$ws.Cells.Item(97, 26).Value = "2025-10-17T07:09:40.120670"
$ws.Cells.Item(98, 26).Value = "2025-10-17T07:09:40.120670"
$ws.Cells.Item(99, 26).Value = "2025-10-17T07:09:40.121211"
$ws.Cells.Item(100, 26).Value = "2025-10-17T07:09:40.121211"
$ws.Cells.Item(101, 26).Value = "2025-10-17T07:09:40.121752"
$ws.Cells.Item(102, 26).Value = "2025-10-17T07:09:40.121752"
$ws.Cells.Item(103, 26).Value = "2025-10-17T07:09:40.175788"
$ws.Cells.Item(104, 26).Value = "2025-10-17T07:09:40.175788"
$ws.Cells.Item(105, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(106, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(107, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(108, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(109, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(110, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(111, 26).Value = "2025-10-17T07:09:40.179803"
$ws.Cells.Item(112, 26).Value = "2025-10-17T07:09:40.179803"
